$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row: "_old" -> "_FV2404" and "_new" -> "_FV2410"
for ($i = 1; $i -le 21; $i++) {
    $cell = $ws.Cells.Item(1, $i)
    $text = $cell.Value2
    if ($text -ne $null) {
        if ($text.EndsWith("_old")) {
            $cell.Value = $text.Substring(0, $text.Length - 4) + "_FV2404"
        } elseif ($text.EndsWith("_new")) {
            $cell.Value = $text.Substring(0, $text.Length - 4) + "_FV2410"
        }
    }
}

# 2. Turn the used range into an Excel Table (ListObject) with default style
$dataRange = $ws.Range("A1:U93")
$listObject = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$listObject.Name = "Table1"
$listObject.TableStyle = ""

# 3. Freeze the header row (split/freeze at row 2)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
